# Auto-generated script applying the scheduled-runner price/profit refresh
# described in the commit "chore: update Sheets via scheduled runner".
# Each hunk in the source diff updates a handful of value cells (H:N) for a
# specific leve row on a specific crafting-class sheet; a few rows also gain
# or lose a LeveProfitNQ/LeveProfitHQ cell entirely.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 3232.0645
$ws.Range("I38").Value = 2300.5
$ws.Range("J38").Value = 3999.2354
$ws.Range("K38").Value = 6901.5
$ws.Range("L38").Value = 11997.7062
$ws.Range("M38").Value = -6529.5
$ws.Range("N38").Value = -12741.7062
$ws.Range("H58").Value = 1131.8572
$ws.Range("I58").Value = 384.6
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 1153.8
$ws.Range("L58").Value = 9000
$ws.Range("M58").Value = -1003.8
$ws.Range("N58").Value = -9300
$ws.Range("H74").Value = 5290.3335
$ws.Range("I74").Value = 4348.6665
$ws.Range("J74").Value = 5761.1665
$ws.Range("K74").Value = 4348.6665
$ws.Range("L74").Value = 5761.1665
$ws.Range("M74").Value = -3412.6665
$ws.Range("N74").Value = -7633.1665
$ws.Range("H77").Value = 5290.3335
$ws.Range("I77").Value = 4348.6665
$ws.Range("J77").Value = 5761.1665
$ws.Range("K77").Value = 21743.3325
$ws.Range("L77").Value = 28805.8325
$ws.Range("M77").Value = -17063.3325
$ws.Range("N77").Value = -38165.8325
$ws.Range("H118").Value = 533.5
$ws.Range("I118").Value = 552.5714
$ws.Range("J118").Value = 400
$ws.Range("K118").Value = 1657.7142
$ws.Range("L118").Value = 1200
$ws.Range("M118").Value = -0.714200000000119
$ws.Range("N118").Value = -4514
$ws.Range("H137").Value = 15649.667
$ws.Range("I137").Value = 20809.412
$ws.Range("J137").Value = 3118.8572
$ws.Range("K137").Value = 62428.236
$ws.Range("L137").Value = 9356.571599999999
$ws.Range("M137").Value = -59878.236
$ws.Range("N137").Value = -14456.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2173.3928
$ws.Range("I74").Value = 2011.6364
$ws.Range("J74").Value = 2766.5
$ws.Range("K74").Value = 2011.6364
$ws.Range("L74").Value = 2766.5
$ws.Range("M74").Value = -1137.6364
$ws.Range("N74").Value = -4514.5
$ws.Range("H77").Value = 2173.3928
$ws.Range("I77").Value = 2011.6364
$ws.Range("J77").Value = 2766.5
$ws.Range("K77").Value = 10058.182
$ws.Range("L77").Value = 13832.5
$ws.Range("M77").Value = -5690.182000000001
$ws.Range("N77").Value = -22568.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 37967.453
$ws.Range("I86").Value = 25649.125
$ws.Range("J86").Value = 70816.336
$ws.Range("K86").Value = 25649.125
$ws.Range("L86").Value = 70816.336
$ws.Range("M86").Value = -24526.125
$ws.Range("N86").Value = -73062.336
$ws.Range("H88").Value = 35000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 35000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 35000
$ws.Range("N88").Value = -35812
$ws.Range("H89").Value = 37967.453
$ws.Range("I89").Value = 25649.125
$ws.Range("J89").Value = 70816.336
$ws.Range("K89").Value = 128245.625
$ws.Range("L89").Value = 354081.68
$ws.Range("M89").Value = -122629.625
$ws.Range("N89").Value = -365313.68
$ws.Range("H91").Value = 35000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 35000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 35000
$ws.Range("N91").Value = -37808

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 1958.8334
$ws.Range("I8").Value = 417.66666
$ws.Range("J8").Value = 3500
$ws.Range("K8").Value = 417.66666
$ws.Range("L8").Value = 3500
$ws.Range("M8").Value = -277.66666
$ws.Range("N8").Value = -3780
$ws.Range("H12").Value = 2998.75
$ws.Range("I12").Value = 397.5
$ws.Range("J12").Value = 5600
$ws.Range("K12").Value = 397.5
$ws.Range("L12").Value = 5600
$ws.Range("M12").Value = -227.5
$ws.Range("N12").Value = -5940
$ws.Range("H13").Value = 5333
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 5333
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 5333
$ws.Range("N13").Value = -5611
$ws.Range("H19").Value = 259.44446
$ws.Range("I19").Value = 227
$ws.Range("J19").Value = 300
$ws.Range("K19").Value = 227
$ws.Range("L19").Value = 300
$ws.Range("M19").Value = -57
$ws.Range("N19").Value = -640
$ws.Range("H24").Value = 259.44446
$ws.Range("I24").Value = 227
$ws.Range("J24").Value = 300
$ws.Range("K24").Value = 227
$ws.Range("L24").Value = 300
$ws.Range("M24").Value = -57
$ws.Range("N24").Value = -640
$ws.Range("H62").Value = 3416.75
$ws.Range("I62").Value = 3359
$ws.Range("J62").Value = 3458
$ws.Range("K62").Value = 3359
$ws.Range("L62").Value = 3458
$ws.Range("M62").Value = -2735
$ws.Range("N62").Value = -4706
$ws.Range("H65").Value = 3416.75
$ws.Range("I65").Value = 3359
$ws.Range("J65").Value = 3458
$ws.Range("K65").Value = 16795
$ws.Range("L65").Value = 17290
$ws.Range("M65").Value = -13675
$ws.Range("N65").Value = -23530
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H141").Value = 82781.10000000001
$ws.Range("I141").Value = 30000
$ws.Range("J141").Value = 85559.05499999999
$ws.Range("K141").Value = 30000
$ws.Range("L141").Value = 85559.05499999999
$ws.Range("M141").Value = -24820
$ws.Range("N141").Value = -95919.05499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 545.6667
$ws.Range("I26").Value = 210
$ws.Range("J26").Value = 713.5
$ws.Range("K26").Value = 630
$ws.Range("L26").Value = 2140.5
$ws.Range("M26").Value = -342
$ws.Range("N26").Value = -2716.5
$ws.Range("H36").Value = 3042.4285
$ws.Range("I36").Value = 998.5
$ws.Range("J36").Value = 3860
$ws.Range("K36").Value = 2995.5
$ws.Range("L36").Value = 11580
$ws.Range("M36").Value = -2826.5
$ws.Range("N36").Value = -11918
$ws.Range("H122").Value = 2028.9231
$ws.Range("I122").Value = 701.5
$ws.Range("J122").Value = 2270.2727
$ws.Range("K122").Value = 6313.5
$ws.Range("L122").Value = 20432.4543
$ws.Range("M122").Value = -3863.5
$ws.Range("N122").Value = -25332.4543

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4848
$ws.Range("I80").Value = 3241.125
$ws.Range("J80").Value = 6454.875
$ws.Range("K80").Value = 3241.125
$ws.Range("L80").Value = 6454.875
$ws.Range("M80").Value = -2243.125
$ws.Range("N80").Value = -8450.875
$ws.Range("H83").Value = 4848
$ws.Range("I83").Value = 3241.125
$ws.Range("J83").Value = 6454.875
$ws.Range("K83").Value = 16205.625
$ws.Range("L83").Value = 32274.375
$ws.Range("M83").Value = -11213.625
$ws.Range("N83").Value = -42258.375
$ws.Range("H113").Value = 1958.8889
$ws.Range("I113").Value = 1763.5834
$ws.Range("J113").Value = 2349.5
$ws.Range("K113").Value = 1763.5834
$ws.Range("L113").Value = 2349.5
$ws.Range("M113").Value = 406.4166
$ws.Range("N113").Value = -6689.5
$ws.Range("H132").Value = 1501.5
$ws.Range("I132").Value = 1273.091
$ws.Range("J132").Value = 4014
$ws.Range("K132").Value = 3819.273
$ws.Range("L132").Value = 12042
$ws.Range("M132").Value = -1289.273
$ws.Range("N132").Value = -17102

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 633969.9399999999
$ws.Range("I93").Value = 1654.4
$ws.Range("J93").Value = 1687829.1
$ws.Range("K93").Value = 1654.4
$ws.Range("L93").Value = 1687829.1
$ws.Range("M93").Value = -406.4000000000001
$ws.Range("N93").Value = -1690325.1
$ws.Range("H132").Value = 2164.303
$ws.Range("I132").Value = 1479.3572
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 4438.071599999999
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -1908.071599999999
$ws.Range("N132").Value = -23060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 101666.25
$ws.Range("I2").Value = 101666.25
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 101666.25
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -101554.25
$ws.Range("H132").Value = 2268.9092
$ws.Range("I132").Value = 2124.9678
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 6374.903399999999
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -3844.903399999999
$ws.Range("N132").Value = -18560
